$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Terminal La Palmera de La Serena - Haba".
# It belongs chronologically before the existing row 54, so insert a fresh row
# there, which pushes the old rows 54-58 down to 55-59 (exactly matching the diff).
$ws.Rows("54:54").Insert()

$ws.Range("A54").Value = 8
$ws.Range("B54").Value = "Terminal La Palmera de La Serena"
$ws.Range("C54").Value = "Coquimbo"
$ws.Range("D54").Value = 45212
$ws.Range("E54").Value = 4
$ws.Range("F54").Value = 100112026
$ws.Range("G54").Value = "Haba"
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 400
$ws.Range("K54").Value = 9000
$ws.Range("L54").Value = 10000
$ws.Range("M54").Value = 9500
$ws.Range("N54").Value = "`$/saco 25 kilos"
$ws.Range("O54").Value = "Provincia del Elquí"
$ws.Range("P54").Value = 380
$ws.Range("Q54").Value = 25
$ws.Range("R54").Value = "Hortaliza"
